$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 216
$ws.Range("E4").Value = 504
$ws.Range("E5").Value = 0
$ws.Range("E6").Value = 504
$ws.Range("E7").Value = 216
$ws.Range("E8").Value = 288
$ws.Range("E9").Value = 144

$ws.Range("E11").Value = 25920
$ws.Range("E12").Value = 60480
$ws.Range("E13").Value = 0
$ws.Range("E14").Value = 60480
$ws.Range("E15").Value = 25920
$ws.Range("E16").Value = 34560
$ws.Range("E17").Value = 17280

$ws.Columns("F:F").Select()
